$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(8)
$ws.Range("A1").Value = "test"
